$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column G: "All" total speed datum
$ws.Range("G1").Value = "All"
$ws.Range("G2").Value = 5
$ws.Range("G3").Formula = "=1/G2*1000"

# Match the number formatting used by the rest of row 3 (integer display)
$ws.Range("G3").NumberFormat = "0"

# Move the active selection, matching the authored workbook state
[void]$ws.Range("J5").Select()
